# Notes_ApresRat.xlsx - rattrapage notes sheet touch-up
#
# - swap the "LastName"/"FirstName" header labels in B1/C1
# - renumber the student CNE codes in column A (rows 2-11): 19000031..19000040
#   become 19000001..19000010 (each shifted down by 30)
# - the renumbered/relabelled cells (A1:C11) pick up a fresh "Normal" style,
#   so make sure they carry the default (un-tinted) formatting
# - leave the cursor on I6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- swap the B1 / C1 header text -----------------------------------------
$b1 = $ws.Range("B1").Value()
$c1 = $ws.Range("C1").Value()
$ws.Range("B1").Value = $c1
$ws.Range("C1").Value = $b1

# --- shift the CNE numbers in column A (rows 2-11) down by 30 -------------
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value() - 30
}

# --- reset A1:C11 to the plain/default cell style --------------------------
$ws.Range("A1:C11").ClearFormats()

# --- move the selection to I6 ----------------------------------------------
$ws.Range("I6").Select()

Write-Output "edit applied"
